$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new data rows (58-110) ---
$ws.Cells.Item(58, 1).Value = 42489
$ws.Cells.Item(58, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(58, 2).Value = 0.54375000000000007
$ws.Cells.Item(58, 2).NumberFormat = "h:mm"
$ws.Cells.Item(58, 3).Value = "Doughnuts"
$ws.Cells.Item(58, 4).Value = $false

$ws.Cells.Item(59, 1).Value = 42487
$ws.Cells.Item(59, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(59, 2).Value = 0.57152777777777775
$ws.Cells.Item(59, 2).NumberFormat = "h:mm"
$ws.Cells.Item(59, 3).Value = "Cake"
$ws.Cells.Item(59, 4).Value = $true

$ws.Cells.Item(60, 1).Value = 42487
$ws.Cells.Item(60, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(60, 2).Value = 0.4145833333333333
$ws.Cells.Item(60, 2).NumberFormat = "h:mm"
$ws.Cells.Item(60, 3).Value = "Cake"
$ws.Cells.Item(60, 4).Value = $true

$ws.Cells.Item(61, 1).Value = 42486
$ws.Cells.Item(61, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(61, 2).Value = 0.54999999999999993
$ws.Cells.Item(61, 2).NumberFormat = "h:mm"
$ws.Cells.Item(61, 3).Value = "Cake"
$ws.Cells.Item(61, 4).Value = $false

$ws.Cells.Item(62, 1).Value = 42481
$ws.Cells.Item(62, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(62, 2).Value = 0.51944444444444449
$ws.Cells.Item(62, 2).NumberFormat = "h:mm"
$ws.Cells.Item(62, 3).Value = "Cake"
$ws.Cells.Item(62, 4).Value = $false

$ws.Cells.Item(63, 1).Value = 42480
$ws.Cells.Item(63, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(63, 2).Value = 0.5493055555555556
$ws.Cells.Item(63, 2).NumberFormat = "h:mm"
$ws.Cells.Item(63, 3).Value = "Cake"
$ws.Cells.Item(63, 4).Value = $true

$ws.Cells.Item(64, 1).Value = 42480
$ws.Cells.Item(64, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(64, 2).Value = 0.34930555555555554
$ws.Cells.Item(64, 2).NumberFormat = "h:mm"
$ws.Cells.Item(64, 3).Value = "Stroopwafels"
$ws.Cells.Item(64, 4).Value = $false

$ws.Cells.Item(65, 1).Value = 42478
$ws.Cells.Item(65, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(65, 2).Value = 0.42083333333333334
$ws.Cells.Item(65, 2).NumberFormat = "h:mm"
$ws.Cells.Item(65, 3).Value = "Pastries"
$ws.Cells.Item(65, 4).Value = $false

$ws.Cells.Item(66, 1).Value = 42473
$ws.Cells.Item(66, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(66, 2).Value = 0.58333333333333337
$ws.Cells.Item(66, 2).NumberFormat = "h:mm"
$ws.Cells.Item(66, 3).Value = "Cake"
$ws.Cells.Item(66, 4).Value = $false

$ws.Cells.Item(67, 1).Value = 42468
$ws.Cells.Item(67, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(67, 2).Value = 0.60555555555555551
$ws.Cells.Item(67, 2).NumberFormat = "h:mm"
$ws.Cells.Item(67, 3).Value = "Cake"
$ws.Cells.Item(67, 4).Value = $true

$ws.Cells.Item(68, 1).Value = 42468
$ws.Cells.Item(68, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(68, 2).Value = 0.58472222222222225
$ws.Cells.Item(68, 2).NumberFormat = "h:mm"
$ws.Cells.Item(68, 3).Value = "Cake"
$ws.Cells.Item(68, 4).Value = $false

$ws.Cells.Item(69, 1).Value = 42465
$ws.Cells.Item(69, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(69, 2).Value = 0.61875000000000002
$ws.Cells.Item(69, 2).NumberFormat = "h:mm"
$ws.Cells.Item(69, 3).Value = "Brownies"
$ws.Cells.Item(69, 4).Value = $false

$ws.Cells.Item(70, 1).Value = 42458
$ws.Cells.Item(70, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(70, 2).Value = 0.6166666666666667
$ws.Cells.Item(70, 2).NumberFormat = "h:mm"
$ws.Cells.Item(70, 3).Value = "Cake"
$ws.Cells.Item(70, 4).Value = $false

$ws.Cells.Item(71, 1).Value = 42458
$ws.Cells.Item(71, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(71, 2).Value = 0.61527777777777781
$ws.Cells.Item(71, 2).NumberFormat = "h:mm"
$ws.Cells.Item(71, 3).Value = "Cake"
$ws.Cells.Item(71, 4).Value = $false

$ws.Cells.Item(72, 1).Value = 42458
$ws.Cells.Item(72, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(72, 2).Value = 0.45277777777777778
$ws.Cells.Item(72, 2).NumberFormat = "h:mm"
$ws.Cells.Item(72, 3).Value = "Sweets"
$ws.Cells.Item(72, 4).Value = $false

$ws.Cells.Item(73, 1).Value = 42452
$ws.Cells.Item(73, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(73, 2).Value = 0.5625
$ws.Cells.Item(73, 2).NumberFormat = "h:mm"
$ws.Cells.Item(73, 3).Value = "Treats"
$ws.Cells.Item(73, 4).Value = $false

$ws.Cells.Item(74, 1).Value = 42452
$ws.Cells.Item(74, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(74, 2).Value = 0.45555555555555555
$ws.Cells.Item(74, 2).NumberFormat = "h:mm"
$ws.Cells.Item(74, 3).Value = "Cake"
$ws.Cells.Item(74, 4).Value = $false

$ws.Cells.Item(75, 1).Value = 42450
$ws.Cells.Item(75, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(75, 2).Value = 0.56597222222222221
$ws.Cells.Item(75, 2).NumberFormat = "h:mm"
$ws.Cells.Item(75, 3).Value = "Treats"
$ws.Cells.Item(75, 4).Value = $false

$ws.Cells.Item(76, 1).Value = 42446
$ws.Cells.Item(76, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(76, 2).Value = 0.61388888888888882
$ws.Cells.Item(76, 2).NumberFormat = "h:mm"
$ws.Cells.Item(76, 3).Value = "Cake"
$ws.Cells.Item(76, 4).Value = $false

$ws.Cells.Item(77, 1).Value = 42446
$ws.Cells.Item(77, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(77, 2).Value = 0.36041666666666666
$ws.Cells.Item(77, 2).NumberFormat = "h:mm"
$ws.Cells.Item(77, 3).Value = "Cookies"
$ws.Cells.Item(77, 4).Value = $false

$ws.Cells.Item(78, 1).Value = 42443
$ws.Cells.Item(78, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(78, 2).Value = 0.55486111111111114
$ws.Cells.Item(78, 2).NumberFormat = "h:mm"
$ws.Cells.Item(78, 3).Value = "Pie"
$ws.Cells.Item(78, 4).Value = $false

$ws.Cells.Item(79, 1).Value = 42437
$ws.Cells.Item(79, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(79, 2).Value = 0.56736111111111109
$ws.Cells.Item(79, 2).NumberFormat = "h:mm"
$ws.Cells.Item(79, 3).Value = "Chocolate"
$ws.Cells.Item(79, 4).Value = $false

$ws.Cells.Item(80, 1).Value = 42437
$ws.Cells.Item(80, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(80, 2).Value = 0.56597222222222221
$ws.Cells.Item(80, 2).NumberFormat = "h:mm"
$ws.Cells.Item(80, 3).Value = "Cake"
$ws.Cells.Item(80, 4).Value = $false

$ws.Cells.Item(81, 1).Value = 42436
$ws.Cells.Item(81, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(81, 2).Value = 0.62916666666666665
$ws.Cells.Item(81, 2).NumberFormat = "h:mm"
$ws.Cells.Item(81, 3).Value = "Brownies"
$ws.Cells.Item(81, 4).Value = $false

$ws.Cells.Item(82, 1).Value = 42432
$ws.Cells.Item(82, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(82, 2).Value = 0.625
$ws.Cells.Item(82, 2).NumberFormat = "h:mm"
$ws.Cells.Item(82, 3).Value = "Cake"
$ws.Cells.Item(82, 4).Value = $false

$ws.Cells.Item(83, 1).Value = 42431
$ws.Cells.Item(83, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(83, 2).Value = 0.44861111111111113
$ws.Cells.Item(83, 2).NumberFormat = "h:mm"
$ws.Cells.Item(83, 3).Value = "Cheesecake"
$ws.Cells.Item(83, 4).Value = $false

$ws.Cells.Item(84, 1).Value = 42431
$ws.Cells.Item(84, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(84, 2).Value = 0.3527777777777778
$ws.Cells.Item(84, 2).NumberFormat = "h:mm"
$ws.Cells.Item(84, 3).Value = "Cake"
$ws.Cells.Item(84, 4).Value = $false

$ws.Cells.Item(85, 1).Value = 42426
$ws.Cells.Item(85, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(85, 2).Value = 0.43611111111111112
$ws.Cells.Item(85, 2).NumberFormat = "h:mm"
$ws.Cells.Item(85, 3).Value = "Cake"
$ws.Cells.Item(85, 4).Value = $false

$ws.Cells.Item(86, 1).Value = 42425
$ws.Cells.Item(86, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(86, 2).Value = 0.62638888888888888
$ws.Cells.Item(86, 2).NumberFormat = "h:mm"
$ws.Cells.Item(86, 3).Value = "Treats"
$ws.Cells.Item(86, 4).Value = $false

$ws.Cells.Item(87, 1).Value = 42422
$ws.Cells.Item(87, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(87, 2).Value = 0.47569444444444442
$ws.Cells.Item(87, 2).NumberFormat = "h:mm"
$ws.Cells.Item(87, 3).Value = "Sweets"
$ws.Cells.Item(87, 4).Value = $false

$ws.Cells.Item(88, 1).Value = 42419
$ws.Cells.Item(88, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(88, 2).Value = 0.3576388888888889
$ws.Cells.Item(88, 2).NumberFormat = "h:mm"
$ws.Cells.Item(88, 3).Value = "Cheesecake"
$ws.Cells.Item(88, 4).Value = $false

$ws.Cells.Item(89, 1).Value = 42419
$ws.Cells.Item(89, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(89, 2).Value = 0.5625
$ws.Cells.Item(89, 2).NumberFormat = "h:mm"
$ws.Cells.Item(89, 3).Value = "Doughnuts"
$ws.Cells.Item(89, 4).Value = $false

$ws.Cells.Item(90, 1).Value = 42419
$ws.Cells.Item(90, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(90, 2).Value = 0.61319444444444449
$ws.Cells.Item(90, 2).NumberFormat = "h:mm"
$ws.Cells.Item(90, 3).Value = "Cookies"
$ws.Cells.Item(90, 4).Value = $false

$ws.Cells.Item(91, 1).Value = 42419
$ws.Cells.Item(91, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(91, 2).Value = 0.6333333333333333
$ws.Cells.Item(91, 2).NumberFormat = "h:mm"
$ws.Cells.Item(91, 3).Value = "Doughnuts"
$ws.Cells.Item(91, 4).Value = $true

$ws.Cells.Item(92, 1).Value = 42411
$ws.Cells.Item(92, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(92, 2).Value = 0.41388888888888892
$ws.Cells.Item(92, 2).NumberFormat = "h:mm"
$ws.Cells.Item(92, 3).Value = "Sweets"
$ws.Cells.Item(92, 4).Value = $false

$ws.Cells.Item(93, 1).Value = 42410
$ws.Cells.Item(93, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(93, 2).Value = 0.6166666666666667
$ws.Cells.Item(93, 2).NumberFormat = "h:mm"
$ws.Cells.Item(93, 3).Value = "Cake"
$ws.Cells.Item(93, 4).Value = $true

$ws.Cells.Item(94, 1).Value = 42405
$ws.Cells.Item(94, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(94, 2).Value = 0.62083333333333335
$ws.Cells.Item(94, 2).NumberFormat = "h:mm"
$ws.Cells.Item(94, 3).Value = "Cake"
$ws.Cells.Item(94, 4).Value = $false

$ws.Cells.Item(95, 1).Value = 42404
$ws.Cells.Item(95, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(95, 2).Value = 0.65486111111111112
$ws.Cells.Item(95, 2).NumberFormat = "h:mm"
$ws.Cells.Item(95, 3).Value = "Cake"
$ws.Cells.Item(95, 4).Value = $true

$ws.Cells.Item(96, 1).Value = 42402
$ws.Cells.Item(96, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(96, 2).Value = 0.4694444444444445
$ws.Cells.Item(96, 2).NumberFormat = "h:mm"
$ws.Cells.Item(96, 3).Value = "Cake"
$ws.Cells.Item(96, 4).Value = $false

$ws.Cells.Item(97, 1).Value = 42398
$ws.Cells.Item(97, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(97, 2).Value = 0.6118055555555556
$ws.Cells.Item(97, 2).NumberFormat = "h:mm"
$ws.Cells.Item(97, 3).Value = "Doughnuts"
$ws.Cells.Item(97, 4).Value = $false

$ws.Cells.Item(98, 1).Value = 42398
$ws.Cells.Item(98, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(98, 2).Value = 0.57013888888888886
$ws.Cells.Item(98, 2).NumberFormat = "h:mm"
$ws.Cells.Item(98, 3).Value = "Cake"
$ws.Cells.Item(98, 4).Value = $true

$ws.Cells.Item(99, 1).Value = 42396
$ws.Cells.Item(99, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(99, 2).Value = 0.68263888888888891
$ws.Cells.Item(99, 2).NumberFormat = "h:mm"
$ws.Cells.Item(99, 3).Value = "Cake"
$ws.Cells.Item(99, 4).Value = $false

$ws.Cells.Item(100, 1).Value = 42391
$ws.Cells.Item(100, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(100, 2).Value = 0.63194444444444442
$ws.Cells.Item(100, 2).NumberFormat = "h:mm"
$ws.Cells.Item(100, 3).Value = "Cake"
$ws.Cells.Item(100, 4).Value = $false

$ws.Cells.Item(101, 1).Value = 42390
$ws.Cells.Item(101, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(101, 2).Value = 0.6743055555555556
$ws.Cells.Item(101, 2).NumberFormat = "h:mm"
$ws.Cells.Item(101, 3).Value = "Cake"
$ws.Cells.Item(101, 4).Value = $false

$ws.Cells.Item(102, 1).Value = 42390
$ws.Cells.Item(102, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(102, 2).Value = 0.63958333333333328
$ws.Cells.Item(102, 2).NumberFormat = "h:mm"
$ws.Cells.Item(102, 3).Value = "Cake"
$ws.Cells.Item(102, 4).Value = $false

$ws.Cells.Item(103, 1).Value = 42390
$ws.Cells.Item(103, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(103, 2).Value = 0.60416666666666663
$ws.Cells.Item(103, 2).NumberFormat = "h:mm"
$ws.Cells.Item(103, 3).Value = "Cake"
$ws.Cells.Item(103, 4).Value = $false

$ws.Cells.Item(104, 1).Value = 42389
$ws.Cells.Item(104, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(104, 2).Value = 0.4368055555555555
$ws.Cells.Item(104, 2).NumberFormat = "h:mm"
$ws.Cells.Item(104, 3).Value = "Cake"
$ws.Cells.Item(104, 4).Value = $false

$ws.Cells.Item(105, 1).Value = 42388
$ws.Cells.Item(105, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(105, 2).Value = 0.60763888888888895
$ws.Cells.Item(105, 2).NumberFormat = "h:mm"
$ws.Cells.Item(105, 3).Value = "Cake"
$ws.Cells.Item(105, 4).Value = $false

$ws.Cells.Item(106, 1).Value = 42387
$ws.Cells.Item(106, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(106, 2).Value = 0.6020833333333333
$ws.Cells.Item(106, 2).NumberFormat = "h:mm"
$ws.Cells.Item(106, 3).Value = "Cake"
$ws.Cells.Item(106, 4).Value = $false

$ws.Cells.Item(107, 1).Value = 42383
$ws.Cells.Item(107, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(107, 2).Value = 0.57430555555555551
$ws.Cells.Item(107, 2).NumberFormat = "h:mm"
$ws.Cells.Item(107, 3).Value = "Cake"
$ws.Cells.Item(107, 4).Value = $true

$ws.Cells.Item(108, 1).Value = 42382
$ws.Cells.Item(108, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(108, 2).Value = 0.40763888888888888
$ws.Cells.Item(108, 2).NumberFormat = "h:mm"
$ws.Cells.Item(108, 3).Value = "Cake"
$ws.Cells.Item(108, 4).Value = $false

$ws.Cells.Item(109, 1).Value = 42376
$ws.Cells.Item(109, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(109, 2).Value = 0.37013888888888885
$ws.Cells.Item(109, 2).NumberFormat = "h:mm"
$ws.Cells.Item(109, 3).Value = "Sweets"
$ws.Cells.Item(109, 4).Value = $false

$ws.Cells.Item(110, 1).Value = 42375
$ws.Cells.Item(110, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(110, 2).Value = 0.42222222222222222
$ws.Cells.Item(110, 2).NumberFormat = "h:mm"
$ws.Cells.Item(110, 3).Value = "Cake"
$ws.Cells.Item(110, 4).Value = $false

# --- Append trailing blank rows (111-134) with date-style formatting on column A ---
for ($r = 111; $r -le 134; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd;@"
}

# --- Correct earlier "Cake" entries that were actually "Treats" ---
$rowsToFix = @(9, 13, 21, 39, 40, 54)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 3).Value = "Treats"
}

# --- Update sheet view: scrolled position and selection ---
$ws.Application.ActiveWindow.ScrollRow = 98
$ws.Range("E110").Select()